$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.270.78"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.884.99"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.00"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.688"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.63"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.353"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.12"
$ws.Range("E10").Value = "  +8.05%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.75"
$ws.Range("E13").Value = "  +7.15%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.778"
$ws.Range("E14").Value = "  +9.89%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.157.73"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.97"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "1.878.56"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "35.266.79"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.16"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.92"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.75"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").Value = "  +6.72%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.01"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.23"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0592"
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.27"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("E33").Value = "  +23.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.15"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -13.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.843"
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0717"
$ws.Range("E39").Value = "  +7.81%  "
$ws.Range("E40").Value = "  +4.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.55"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.09"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "1.325.64"
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.20"
$ws.Range("E45").Value = "  +12.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0806"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "2.057.36"
$ws.Range("E51").Value = "  +0.04%  "
